# Työaikakirjanpito edit script
# Updates time-tracking log rows: fixes a few earlier entries' text/hours,
# and fills in three new rows (9-11) of logged work, per commit
# "validaattori testit etc Github sekos".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: "Firebasen opettelua ja mapsAPIn toiminnan selvittämistä" -> "Firebasen opettelua"
$ws.Cells.Item(3, 3).Value = "Firebasen opettelua"

# Row 4: Three.js text updated
$ws.Cells.Item(4, 3).Value = "Three.js:n harjoittelua ja demoamista. Mahdollisesti animaatioiden tekoa sivulle"

# Row 6: hours 7.5 -> 3.5, text "Blender objectien..." -> "Three.js:n objectien muodostamista"
$ws.Cells.Item(6, 2).Value = 3.5
$ws.Cells.Item(6, 3).Value = "Three.js:n objectien muodostamista"

# Row 9: new entry
$ws.Cells.Item(9, 1).Value = 44025
$ws.Cells.Item(9, 2).Value = 1.5
$ws.Cells.Item(9, 3).Value = "Firebase with react"

# Row 10: new entry
$ws.Cells.Item(10, 1).Value = 44027
$ws.Cells.Item(10, 2).Value = 6
$ws.Cells.Item(10, 3).Value = "3d mallien importtaaminen ja menun rakentaminen"

# Row 11: new entry
$ws.Cells.Item(11, 1).Value = 44033
$ws.Cells.Item(11, 2).Value = 4.5
$ws.Cells.Item(11, 3).Value = "Databasen rakentamista, validaattoreidenn tekemistä"

# Row heights: row 3 shrinks back to the default (autofit clears the
# explicit height), row 4 shrinks from 45 to 30, rows 10-11 (now with
# wrapped text) grow to 30.
$ws.Rows.Item(3).AutoFit()
$ws.Rows.Item(4).RowHeight = 30
$ws.Rows.Item(10).RowHeight = 30
$ws.Rows.Item(11).RowHeight = 30

# Selection moved to F10 in the saved view.
[void]$ws.Range("F10").Select()
